$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("","","","","","")
    ,@("20104155","NISSIN HOT CARBO 120","EG1B","1","1","RT,(E-1B)")
    ,@("20101747","INDOMI SBLAK HOT/J75","EG1B","1","2","RT,(E-2B)")
    ,@("20138358","INDOMI NYEMK B/LDS80","EG1B","1","3","RT,(E-2B)")
    ,@("10007606","ABC MIE GULAI PDS.70","EG1B","1","4","RT,(E-1B)")
    ,@("10022534","ABC KCP MNS PCH685G","EG1B","1","5","RT,(E-2B)")
    ,@("20054589","3SAPI KRIMER 490G","EG1B","1","6","RT,(E-1B)")
    ,@("20072194","MH MAITOS S.BLDO 140","EG1B","2","1","RT,(E-3B)")
    ,@("20116675","CHTATO LITE NORI 115","EG1B","2","2","RT,(E-14H)")
    ,@("20116112","MAXICORN BBQ 140G","EG1B","2","3","RT,(E-14H)")
    ,@("20104728","CHIMI KRP.UB JG.BK41","EG1B","2","4","RT,(E-3B)")
    ,@("20138228","LAYS RUMPUT LAUT 64G","EG1B","3","1","RT,(E-1B)")
    ,@("20114756","GURIBEE BBQ BLDO 65G","EG1B","3","2","RT,(E-14H)")
    ,@("10000753","TARO SNACK SEAWED 32","EG1B","3","3","RT,(E-2B)")
    ,@("20045104","QTELA TEMPE CABE/R55","EG1B","3","4","RT,(E-14H)")
    ,@("20124895","G/T COKIES CHOCDIP71","EG1B","4","1","RT,(E-2B)")
    ,@("10035165","OREO ORIGINAL 105G","EG1B","4","2","RT,(E-1B)")
    ,@("20057475","SERENA TOGO CKLT 128","EG1B","4","3","RT,(E-1B)")
    ,@("20048791","OREO MINI VANLA 58.4","EG1B","4","4","RT,(E-1B)")
    ,@("20098348","REBO KUACI ORIGNL120","EG1B","5","1","RT,(E-1B)")
    ,@("20093909","REBO KUACI MILK 120","EG1B","5","2","RT,(E-1B)")
    ,@("20133506","DELFI KRICE SWET64.8","EG1B","5","3","RT,(E-1B)")
    ,@("20133507","DELFI KRICE SWD 47.4","EG1B","5","4","RT,(E-1B)")
    ,@("20065180","GLICO PEJOY CHO 30G","EG1B","5","5","RT,(E-1B)")
    ,@("20135539","MR.POTATO ORGINAL85G","EG1B","5","6","RT,(E-1B)")
    ,@("20136007","MR.POTATO BBQ 85G","EG1B","5","7","RT,(E-1B)")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# Sheet shrank from 27 rows (A1:F27) to 26 rows (A1:F26); remove the leftover last row
$ws.Rows.Item(27).Delete()

# Column F width changed from 13 to 12
$ws.Columns.Item(6).ColumnWidth = 12
